$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.16594766666667
$ws.Range("H2").Value = 39.497843
$ws.Range("I2").Value = 0.6940777873489595
$ws.Range("J2").Value = 0.6940777873489595
$ws.Range("M2").Value = 0.06743766666666666
$ws.Range("N2").Value = 0.202313
$ws.Range("O2").Value = 0.004349811883262162
$ws.Range("P2").Value = 0.004349811883262163
$ws.Range("Q2").Value = 0.8878807900954444
$ws.Range("R2").Value = 7.990927110859
$ws.Range("S2").Value = 0.003019107807318813
$ws.Range("T2").Value = 0.003019107807318813
$ws.Range("G3").Value = 13.16594766666667
$ws.Range("H3").Value = 39.497843
$ws.Range("I3").Value = 0.6940777873489595
$ws.Range("J3").Value = 0.6940777873489595
$ws.Range("O3").Value = 0.7938207485680675
$ws.Range("P3").Value = 0.7938207485680676
$ws.Range("Q3").Value = 162.0341781089142
$ws.Range("R3").Value = 1458.307602980228
$ws.Range("S3").Value = 0.550973348717819
$ws.Range("T3").Value = 0.5509733487178191
$ws.Range("G4").Value = 13.16594766666667
$ws.Range("H4").Value = 39.497843
$ws.Range("I4").Value = 0.6940777873489595
$ws.Range("J4").Value = 0.6940777873489595
$ws.Range("O4").Value = 0.2018294395486703
$ws.Range("P4").Value = 0.2018294395486704
$ws.Range("Q4").Value = 41.19729474751489
$ws.Range("R4").Value = 370.7756527276341
$ws.Range("S4").Value = 0.1400853308238217
$ws.Range("T4").Value = 0.1400853308238217
$ws.Range("I5").Value = 0.1706596770095176
$ws.Range("J5").Value = 0.1706596770095176
$ws.Range("M5").Value = 0.06743766666666666
$ws.Range("N5").Value = 0.202313
$ws.Range("O5").Value = 0.004349811883262162
$ws.Range("P5").Value = 0.004349811883262163
$ws.Range("Q5").Value = 0.21831191204
$ws.Range("R5").Value = 1.96480720836
$ws.Range("S5").Value = 0.0007423374910496822
$ws.Range("T5").Value = 0.0007423374910496824
$ws.Range("I6").Value = 0.1706596770095176
$ws.Range("J6").Value = 0.1706596770095176
$ws.Range("O6").Value = 0.7938207485680675
$ws.Range("P6").Value = 0.7938207485680676
$ws.Range("S6").Value = 0.1354731925540799
$ws.Range("T6").Value = 0.1354731925540799
$ws.Range("I7").Value = 0.1706596770095176
$ws.Range("J7").Value = 0.1706596770095176
$ws.Range("O7").Value = 0.2018294395486703
$ws.Range("P7").Value = 0.2018294395486704
$ws.Range("S7").Value = 0.03444414696438804
$ws.Range("T7").Value = 0.03444414696438804
$ws.Range("H8").Value = 7.697376999999999
$ws.Range("I8").Value = 0.1352625356415228
$ws.Range("J8").Value = 0.1352625356415228
$ws.Range("M8").Value = 0.06743766666666666
$ws.Range("N8").Value = 0.202313
$ws.Range("O8").Value = 0.004349811883262162
$ws.Range("P8").Value = 0.004349811883262163
$ws.Range("Q8").Value = 0.1730310481112222
$ws.Range("R8").Value = 1.557279433001
$ws.Range("S8").Value = 0.0005883665848936676
$ws.Range("T8").Value = 0.0005883665848936678
$ws.Range("H9").Value = 7.697376999999999
$ws.Range("I9").Value = 0.1352625356415228
$ws.Range("J9").Value = 0.1352625356415228
$ws.Range("O9").Value = 0.7938207485680675
$ws.Range("P9").Value = 0.7938207485680676
$ws.Range("S9").Value = 0.1073742072961685
$ws.Range("T9").Value = 0.1073742072961685
$ws.Range("H10").Value = 7.697376999999999
$ws.Range("I10").Value = 0.1352625356415228
$ws.Range("J10").Value = 0.1352625356415228
$ws.Range("O10").Value = 0.2018294395486703
$ws.Range("P10").Value = 0.2018294395486704
$ws.Range("Q10").Value = 8.028567763858444
$ws.Range("S10").Value = 0.02729996176046058
$ws.Range("T10").Value = 0.02729996176046059
